$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 60: Thursday 2024-07-25, 2 hours
$ws.Range("A60").Value = 45498
$ws.Range("B60").Value = "T"
$ws.Range("C60").Value = 2
$ws.Range("E60").Value = "Cleaning up graph comparing models, sent email, reading up on imputation"

# Row 61: Friday 2024-07-26, 2 hours
$ws.Range("A61").Value = 45499
$ws.Range("B61").Value = "F"
$ws.Range("C61").Value = 2
$ws.Range("E61").Value = "Testing studentgrades_prof with imputation again, on all courses used in GBM"

# Match formatting used by the rest of the table for new rows
$ws.Range("A60").NumberFormat = $ws.Range("A59").NumberFormat
$ws.Range("A61").NumberFormat = $ws.Range("A59").NumberFormat
$ws.Range("B60").NumberFormat = $ws.Range("B59").NumberFormat
$ws.Range("B61").NumberFormat = $ws.Range("B59").NumberFormat
$ws.Range("E60").WrapText = $true
$ws.Range("E61").WrapText = $true

$ws.Rows.Item(60).RowHeight = 30
$ws.Rows.Item(61).RowHeight = 30

$ws.Range("F59").Select()
